$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of data (row 6) mirroring the existing rows, using the
# same formatting already applied to the date column (A) on row 5 so
# that the new cell reuses the existing style record instead of creating
# a new one.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A6").Value = 42607.889236111114

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 60
$ws.Range("D6").Value = 38
$ws.Range("E6").Value = 31
$ws.Range("F6").Value = 68
$ws.Range("G6").Value = 42328
$ws.Range("H6").Value = 20478
$ws.Range("I6").Value = 3780
$ws.Range("J6").Value = 381
$ws.Range("K6").Value = 240
$ws.Range("L6").Value = 14
$ws.Range("M6").Value = 31
$ws.Range("N6").Value = "Noun"
